$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New validation columns (X, Y) — header row first, then data row values in the
# exact order needed so the shared-string table grows in the same sequence as
# the authored edit.
$ws.Range("X1").Value = "validProgId_User"
$ws.Range("X2").Value = "'16265"
$ws.Range("Y2").Value = "'8470"
$ws.Range("Y1").Value = "validBatchId_User"

# Updated data-row values for the existing columns.
$ws.Range("Q2").Value = "R02"
$ws.Range("D2").Value = "SDET-007"
$ws.Range("K2").Value = "'LMSstaff"
$ws.Range("L2").Value = "'NumpyNinja"
$ws.Range("O2").Value = "Teamfourteen"
$ws.Range("A2").Value = "Team14-REST-riders-SDET-51"
$ws.Range("N2").Value = "restriders_team14_Staff_j@gmail.com"
$ws.Range("P2").Value = "'5754525409"

# Column widths were nudged by Excel as the new data/columns were entered.
$ws.Columns.Item(16).ColumnWidth = 14.92
$ws.Columns.Item(17).ColumnWidth = 7.42
$ws.Columns.Item(18).ColumnWidth = 10.1
$ws.Columns.Item(24).ColumnWidth = 15.92
$ws.Columns.Item(25).ColumnWidth = 18.25

# Match the author's final selection/navigation state (scrolled right, cell
# P5 active) — topLeftCell scroll position isn't exposed/persisted by this
# host, so only the active-cell selection is reproducible here.
$ws.Range("P5").Select()
